# Update device_type.xlsx for release 1.1.5:
# Add French (fra) translations for each device type, interleaved with the
# existing English (eng) rows, so each device code has an eng row followed
# by a fra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: French translation of FRS (Finger Print Scanner)
$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "FRS"
$ws.Range("C3").Value = "Scanner dempreintes digitales"
$ws.Range("D3").Value = "Scannez les empreintes digitales"
# E3 already contains "TRUE" - leave untouched

# Row 4: English IRS (Iris Scanner) - shifted down from its old row 3 slot
$ws.Range("A4").Value = "eng"
$ws.Range("B4").Value = "IRS"
$ws.Range("C4").Value = "Iris Scanner"
$ws.Range("D4").Value = "For scanning Iris"
# E4 already contains "TRUE" - leave untouched

# Row 5: French translation of IRS (Iris Scanner)
$ws.Range("A5").Value = "fra"
$ws.Range("B5").Value = "IRS"
$ws.Range("C5").Value = "Scanner dIris"
$ws.Range("D5").Value = "Pour scanner liris"
# E5 already contains "TRUE" - leave untouched

# Row 6: English CMR (Camera) - shifted down from its old row 4 slot
$ws.Range("A6").Value = "eng"
$ws.Range("B6").Value = "CMR"
$ws.Range("C6").Value = "Camera"
$ws.Range("D6").Value = "For capturing photo"
# E6 already contains "TRUE" - leave untouched

# Row 7: French translation of CMR (Camera)
$ws.Range("A7").Value = "fra"
$ws.Range("B7").Value = "CMR"
$ws.Range("C7").Value = "Caméra"
$ws.Range("D7").Value = "Pour capturer une photo"
$ws.Range("E7").Value = "TRUE"

# Row 8: English SCN (Document Scanner) - shifted down from its old row 5 slot
$ws.Range("A8").Value = "eng"
$ws.Range("B8").Value = "SCN"
$ws.Range("C8").Value = "Document Scanner"
$ws.Range("D8").Value = "For scanning documents"
$ws.Range("E8").Value = "TRUE"

# Row 9: French translation of SCN (Document Scanner)
$ws.Range("A9").Value = "fra"
$ws.Range("B9").Value = "SCN"
$ws.Range("C9").Value = "Scanner de documents"
$ws.Range("D9").Value = "Pour numériser des documents"
$ws.Range("E9").Value = "TRUE"

# Row 10: English PRT (Printer) - shifted down from its old row 6 slot
$ws.Range("A10").Value = "eng"
$ws.Range("B10").Value = "PRT"
$ws.Range("C10").Value = "Printer"
$ws.Range("D10").Value = "For printing Documents"
$ws.Range("E10").Value = "TRUE"

# Row 11: French translation of PRT (Printer)
$ws.Range("A11").Value = "fra"
$ws.Range("B11").Value = "PRT"
$ws.Range("C11").Value = "Imprimante"
$ws.Range("D11").Value = "Pour imprimer des documents"
$ws.Range("E11").Value = "TRUE"

$ws.Range("G10").Select()
